# Financials update: adjust Balance Sheet figures on the CMTA sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Other Current Assets (row 45): first-period value revised
$ws.Range("D45").Value = 3000

# Total Current Assets (row 46): first-period value revised
$ws.Range("D46").Value = 70900

# Goodwill (row 49): first-period value revised
$ws.Range("D49").Value = 3400

# Other Assets (row 52): now has a value in the first period, and the
# remaining periods are marked "NA" instead of 0
$ws.Range("D52").Value = 800
$ws.Range("E52:J52").Value = "NA"

# Capital Expenditures (row 91): zero out the negative figures
$ws.Range("D91").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
